$d = $word.ActiveDocument

# Locate the paragraph that ends the course entry ("LOQ4240: ... (Requisito
# fraco)"). This is the last piece of real course content before the page
# footer block that needs to be stripped.
$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "LOQ4240*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the LOQ4240 paragraph anchor."
}

# Immediately after that paragraph the page used to carry three more
# paragraphs that are now removed from the site footer:
#   1. a blank spacer paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# Delete that whole span (its own blank paragraph + trailing page-break
# paragraph remain untouched).
$firstToDelete = $d.Paragraphs.Item($anchorIndex + 1)
$lastToDelete = $d.Paragraphs.Item($anchorIndex + 3)

$r = $d.Range($firstToDelete.Range.Start, $lastToDelete.Range.End)
$r.Delete()
